# Última atualização dos casos de uso.
#
# 1) Three of the five "Passo 5: " step headings (the ones that lead into
#    the "Excluir", "Alterar" and "Voltar" alternative-flow steps) become
#    "Passo 3: " -- the two "Passo 5" occurrences that read "...Cancelar"
#    and "...da Seção Principal" are left untouched.
# 2) "Deseja realmente excluir a Fornecedor?" becomes "...excluir o
#    Fornecedor?" (article "a" -> "o").
#
# Word materialises each of these single-character substitutions as its
# own run (identical rPr to its neighbours) rather than re-merging back
# into the surrounding run, so we reproduce that by toggling Bold off/on
# immediately after the substitution -- this forces the run boundary
# without leaving any residual direct formatting behind.

$d = $word.ActiveDocument

function Get-MatchStart($needle) {
    $rng = $d.Content
    $rng.Start = 0
    $rng.End = $d.Content.End
    if ($rng.Find.Execute($needle, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)) {
        return $rng.Start
    }
    return -1
}

# --- "Passo 5: " -> "Passo 3: " for the three steps that continue with
#     "Excluir" / "Alterar" / "Voltar".
$afterTexts = @(
    "Ator pressiona o botão “Excluir”. Ver",
    "Ator pressiona o botão “Alterar”. Ver",
    "Ator pressiona botão “Voltar”. Retorna"
)

foreach ($after in $afterTexts) {
    $needle = "Passo 5: " + $after
    $start = Get-MatchStart($needle)
    if ($start -lt 0) {
        Write-Output "NOT FOUND: $needle"
        continue
    }

    # "5" is the 7th character (0-based index 6) of "Passo 5: ".
    $digitRange = $d.Range($start + 6, $start + 7)
    $digitRange.Text = "3"

    # Re-acquire the (now "3") one-character range and toggle Bold off/on
    # so the edited character ends up in its own run, matching how Word
    # itself splits runs on a same-formatted in-place retype.
    $digitRange = $d.Range($start + 6, $start + 7)
    $digitRange.Font.Bold = $false
    $digitRange.Font.Bold = $true
}

# --- "excluir a Fornecedor" -> "excluir o Fornecedor".
$needle2 = "excluir a Fornecedor"
$start2 = Get-MatchStart($needle2)
if ($start2 -ge 0) {
    # "a" is the 9th character (0-based index 8) of "excluir a Fornecedor".
    $letterRange = $d.Range($start2 + 8, $start2 + 9)
    $letterRange.Text = "o"

    # Same trick as above, using Bold (this run isn't bold, so the final
    # state matches the original and leaves no stray formatting behind).
    $letterRange = $d.Range($start2 + 8, $start2 + 9)
    $letterRange.Font.Bold = $true
    $letterRange.Font.Bold = $false
} else {
    Write-Output "NOT FOUND: $needle2"
}
